$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($rangeRef, $val) {
    $rng = $ws.Range($rangeRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue "D2" "38.024.68"
Set-TextValue "E2" "  +2.35%  "
Set-TextValue "D3" "2.054.11"
Set-TextValue "E3" "  +1.41%  "
Set-TextValue "E4" "  +0.13%  "
Set-TextValue "D5" "229.55"
Set-TextValue "E5" "  +1.15%  "
Set-TextValue "D6" "0.616"
Set-TextValue "E6" "  +1.20%  "
Set-TextValue "D7" "58.75"
Set-TextValue "E7" "  +6.68%  "
Set-TextValue "E8" "  -0.05%  "
Set-TextValue "D9" "0.386"
Set-TextValue "E9" "  +1.79%  "
Set-TextValue "D10" "0.0810"
Set-TextValue "E10" "  +2.93%  "
Set-TextValue "E11" "  +1.69%  "
Set-TextValue "D12" "2.356.61"
Set-TextValue "E12" "  +1.72%  "
Set-TextValue "D13" "14.67"
Set-TextValue "E13" "  +2.99%  "
Set-TextValue "D14" "20.86"
Set-TextValue "E14" "  +3.00%  "
Set-TextValue "B15" "Polygon"
Set-TextValue "C15" "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
Set-TextValue "D15" "0.752"
Set-TextValue "E15" "  +1.31%  "
Set-TextValue "B16" "Polkadot"
Set-TextValue "C16" "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextValue "D16" "5.29"
Set-TextValue "E16" "  +1.84%  "
Set-TextValue "D17" "2.054.03"
Set-TextValue "E17" "  +3.71%  "
Set-TextValue "D18" "37.946.18"
Set-TextValue "E18" "  +2.22%  "
Set-TextValue "D19" "6.25"
Set-TextValue "E19" "  -3.20%  "
Set-TextValue "D20" "69.71"
Set-TextValue "E20" "  +1.20%  "
Set-TextValue "E21" "  +2.27%  "
Set-TextValue "D22" "224.71"
Set-TextValue "E22" "  +0.42%  "
Set-TextValue "E23" "  +0.09%  "
Set-TextValue "E24" "  -0.07%  "
Set-TextValue "D25" "2.24"
Set-TextValue "E25" "  +2.41%  "
Set-TextValue "B26" "Cosmos"
Set-TextValue "C26" "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextValue "D26" "9.29"
Set-TextValue "E26" "  +0.93%  "
Set-TextValue "B27" "Monero"
Set-TextValue "C27" "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue "D27" "166.26"
Set-TextValue "E27" "  +0.26%  "
Set-TextValue "D28" "0.134"
Set-TextValue "E28" "  +4.50%  "
Set-TextValue "D29" "19.02"
Set-TextValue "E29" "  +1.54%  "
Set-TextValue "E30" "  +0.30%  "
Set-TextValue "E31" "  +1.42%  "
Set-TextValue "D32" "4.53"
Set-TextValue "E32" "  +0.22%  "
Set-TextValue "D33" "4.60"
Set-TextValue "E33" "  +2.46%  "
Set-TextValue "D34" "2.06"
Set-TextValue "E34" "  +10.59%  "
Set-TextValue "E35" "  -0.63%  "
Set-TextValue "E36" "  -0.68%  "
Set-TextValue "E37" "  +9.96%  "
Set-TextValue "E38" "  +5.36%  "
Set-TextValue "E39" "  +0.01%  "
Set-TextValue "D40" "1.488.13"
Set-TextValue "E40" "  +1.20%  "
Set-TextValue "E41" "  +0.83%  "
Set-TextValue "D42" "97.09"
Set-TextValue "E42" "  +1.32%  "
Set-TextValue "D43" "2.84"
Set-TextValue "E43" "  +0.85%  "
Set-TextValue "D44" "16.56"
Set-TextValue "E44" "  +0.69%  "
Set-TextValue "D45" "0.0924"
Set-TextValue "E45" "  +1.36%  "
Set-TextValue "E46" "  -1.14%  "
Set-TextValue "D47" "4.14"
Set-TextValue "E47" "  +13.99%  "
Set-TextValue "E48" "  +0.51%  "
Set-TextValue "E49" "  +1.16%  "
Set-TextValue "D50" "7.06"
Set-TextValue "E50" "  -2.72%  "
Set-TextValue "D51" "2.243.28"
Set-TextValue "E51" "  +1.48%  "
